# Weekly crime-data refresh for the 60th Precinct CompStat report.
# Updates the report header (volume number + week-covering date range)
# and the full crime-complaints grid (rows 15-31) with the newly
# collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: Volume/Number + reporting week range ---------------------
$ws.Range("A8").Value = "Volume 32   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/16/2025  Through  6/22/2025"

# --- Burglary row (18): "n/a" placeholder -> 1, with numeric styling --
# Pull formatting (style 14, "#,##0") from a sibling numeric cell in the
# same row so the new value renders like the rest of the grid.
$ws.Range("I18").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 1

# --- Transit row (22): 1 -> "n/a" placeholder, with text styling ------
# Pull formatting (style 13, general/text, shares the "0" placeholder
# string) from a sibling cell that already displays the placeholder.
$ws.Range("D27").Copy($ws.Range("C22"))

# --- Crime-complaints grid: weekly/28-day/YTD/2-year counts and the
#     associated percent-change figures ------------------------------
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 20
$ws.Range("K15").Value = 233.333333333333
$ws.Range("L15").Value = 400
$ws.Range("M15").Value = 1900
$ws.Range("N15").Value = 33.333333333333
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 13
$ws.Range("E16").Value = -61.538461538461
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 60
$ws.Range("J16").Value = 77
$ws.Range("K16").Value = -22.077922077922
$ws.Range("L16").Value = -18.918918918918
$ws.Range("M16").Value = -47.368421052631
$ws.Range("N16").Value = -87.730061349693
$ws.Range("C17").Value = 7
$ws.Range("E17").Value = 16.666666666666
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = -16.216216216216
$ws.Range("I17").Value = 190
$ws.Range("J17").Value = 169
$ws.Range("K17").Value = 12.426035502958
$ws.Range("L17").Value = 15.151515151515
$ws.Range("M17").Value = 118.390804597701
$ws.Range("N17").Value = -43.620178041543
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 4
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 52
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = 30
$ws.Range("L18").Value = 15.555555555555
$ws.Range("M18").Value = -20
$ws.Range("N18").Value = -89.278350515463
$ws.Range("C19").Value = 7
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = -26.190476190476
$ws.Range("I19").Value = 143
$ws.Range("J19").Value = 168
$ws.Range("K19").Value = -14.880952380952
$ws.Range("L19").Value = -35.294117647058
$ws.Range("M19").Value = -38.626609442060
$ws.Range("N19").Value = -48.745519713261
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -18.181818181818
$ws.Range("I20").Value = 47
$ws.Range("J20").Value = 53
$ws.Range("K20").Value = -11.320754716981
$ws.Range("L20").Value = 17.5
$ws.Range("M20").Value = -7.843137254901
$ws.Range("N20").Value = -90.637450199203
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = -41.025641025641
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 118
$ws.Range("H21").Value = -22.881355932203
$ws.Range("I21").Value = 515
$ws.Range("J21").Value = 514
$ws.Range("K21").Value = 0.194552529182
$ws.Range("L21").Value = -6.871609403254
$ws.Range("M21").Value = -6.871609403254
$ws.Range("N21").Value = -75.673122342938
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -57.142857142857
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 88.888888888888
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -7.692307692307
$ws.Range("I23").Value = 71
$ws.Range("J23").Value = 71
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1.428571428571
$ws.Range("M23").Value = 57.777777777777
$ws.Range("C24").Value = 20
$ws.Range("E24").Value = -4.761904761904
$ws.Range("F24").Value = 67
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = -19.277108433734
$ws.Range("I24").Value = 433
$ws.Range("J24").Value = 455
$ws.Range("K24").Value = -4.835164835164
$ws.Range("L24").Value = -8.067940552016
$ws.Range("M24").Value = -10.537190082644
$ws.Range("C25").Value = 4
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 18
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = 12.5
$ws.Range("I25").Value = 119
$ws.Range("J25").Value = 123
$ws.Range("K25").Value = -3.252032520325
$ws.Range("L25").Value = -26.993865030674
$ws.Range("C26").Value = 22
$ws.Range("D26").Value = 23
$ws.Range("E26").Value = -4.347826086956
$ws.Range("F26").Value = 68
$ws.Range("G26").Value = 75
$ws.Range("H26").Value = -9.333333333333
$ws.Range("I26").Value = 345
$ws.Range("J26").Value = 352
$ws.Range("K26").Value = -1.988636363636
$ws.Range("L26").Value = 28.731343283582
$ws.Range("M26").Value = 52.654867256637
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 150
$ws.Range("I27").Value = 22
$ws.Range("K27").Value = 69.230769230769
$ws.Range("L27").Value = 175
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 10
$ws.Range("H28").Value = 11.111111111111
$ws.Range("I28").Value = 38
$ws.Range("J28").Value = 27
$ws.Range("K28").Value = 40.740740740740
$ws.Range("L28").Value = 46.153846153846
$ws.Range("L31").Value = -11.111111111111
